$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 22:52"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 528990
$ws.Range("C4").Value = 26114
$ws.Range("D4").Value = 29436
$ws.Range("E4").Value = 479099
$ws.Range("F4").Value = 11057
$ws.Range("G4").Value = 1708
$ws.Range("H4").Value = 20455

# Row 29 - Noruega
$ws.Range("E29").Value = 6257
$ws.Range("G29").Value = 6
$ws.Range("H29").Value = 119

# Row 74 - Bosnia y Herzegovina
$ws.Range("B74").Value = 946
$ws.Range("C74").Value = 45
$ws.Range("E74").Value = 770

# Row 106 - Nigeria
$ws.Range("B106").Value = 318
$ws.Range("C106").Value = 13
$ws.Range("D106").Value = 70
$ws.Range("E106").Value = 238
$ws.Range("G106").Value = 3
$ws.Range("H106").Value = 10
